$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.468.15"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.106.01"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "334.01"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("D7").Value = "0.5223"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.4527"
$ws.Range("E8").Value = "  +4.80%  "
$ws.Range("D9").Value = "53.19"
$ws.Range("E9").Value = "  +14.46%  "
$ws.Range("D10").Value = "0.08916"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "1.178"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").Value = "2.101.95"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "6.826"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").Value = "8.032"
$ws.Range("E15").Value = "  +4.86%  "
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "0.00001143"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "0.06652"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "6.342"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "30.530.41"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "12.46"
$ws.Range("D25").Value = "2.346"
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("D26").Value = "2.351.16"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").Value = "22.23"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").Value = "162.83"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "2.530"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").Value = "133.19"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").Value = "1.207"
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").Value = "1.661"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").Value = "6.428"
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("D35").Value = "3.941"
$ws.Range("E35").Value = "  +2.76%  "
$ws.Range("D36").Value = "10.43"
$ws.Range("E36").Value = "  +5.23%  "
$ws.Range("D37").Value = "5.809"
$ws.Range("E37").Value = "  +6.62%  "
$ws.Range("D38").Value = "0.02594"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").Value = "12.73"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "0.6876"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("D43").Value = "1.250"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "2.312"
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.99"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "0.6368"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").Value = "0.00000000351"
$ws.Range("E48").Value = "  +21.40%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "83.51"
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "1.204"
$ws.Range("E51").Value = "  +1.52%  "
